$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new (blank) rows right above the "DHT 1" block (old row 4 -> new row 6)
$ws.Rows("4:5").Insert()

# Insert 4 new (blank) rows between the two hash/socket rows and the
# "INFO DE REDUNDANCIA POR CAIDA" block (old row 7 -> new row 13)
$ws.Rows("9:12").Insert()

# The first hash/socket row picks up an explicit (default-sized) row height
$ws.Rows("7").RowHeight = 15

# The old D7:D11 merged cell (now shifted to D13:D17) is split apart again
$ws.Range("D13:D17").UnMerge()

# That row now holds the full (previously merged/wrapped) text, so it grows tall
$ws.Rows("13").RowHeight = 60

# New port numbers next to the hash rows
$ws.Range("G7").Value = 64
$ws.Range("G8").Value = 69

# New third hash row
$ws.Range("B9").Value = "HASH 3"
$ws.Range("G9").Value = 64

# New empty, underline-styled placeholder cells (style matches the "HASH n" header cells)
$ws.Range("I10").Font.Underline = $true
$ws.Range("I11").Font.Underline = $true
$ws.Range("H13").Font.Underline = $true

# Match the author's final selection
$ws.Range("I10:I11").Select() | Out-Null
